$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix unit/name issue on DPIE data: correct AED Name for Chla_mg/m3 row
# from "WQ_DIAG_PHY_TCHA" to "WQ_DIAG_PHY_TCHLA"
$ws.Range("B17").Value = "WQ_DIAG_PHY_TCHLA"
